$d = $word.ActiveDocument

# Locate the keywords paragraph that currently reads
# " Reproducción, seres vivos," (spread across two runs in the source).
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" Reproducción, seres vivos,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target keywords text to edit"
}

$matchStart = $rng.Start
$matchEnd = $rng.End

# The paragraph also carries a zero-width "_GoBack" bookmark right at the
# end of the matched text (before the paragraph mark). Remove it first so
# it doesn't get auto-preserved/duplicated when we overwrite the range's
# XML; we'll re-insert it at the correct spot explicitly below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-anchor on fresh character offsets: the previously fetched range can
# go stale once the document is mutated by the bookmark deletion above.
$rng = $d.Range($matchStart, $matchEnd)

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-ES_tradnl"/></w:rPr>'

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>r</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r>' + $rPr + '<w:t>eproducci&#243;n,</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>seres</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> vivos</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xmlFrag)
